$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weapon_Stat")

# Update Index values for rows 3 and 4 (A3, A4)
$ws.Range("A3").Value = 2001
$ws.Range("A4").Value = 2002

# Move the active selection from D5 to D6
$ws.Range("D6").Select()
